$wb = $excel.ActiveWorkbook

$commitBase   = "https://github.com/OpenLocalizationTest/oltest/blob/39d08eb83a2fd194452a4bab940ad0d343190131/"
$e2eBase      = $commitBase + "e2e/"
$handoffZh    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/100a81bb83ebc05a0e466446f84d5cf383345b20/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$handoffDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b1eaa0bc8dee9175e2e461a277ca2797fcc5c0a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$newMdA = "9e58fb96-cb81-4926-9015-4068e973d72d.md"
$newMdB = "a07768e3-e4d7-49bf-90a4-2429cf1f557e.md"
$xlfZhA = "9e58fb96-cb81-4926-9015-4068e973d72d.81b0d9db6138142e445d0a5a643a729d74d95ba8.zh-cn.xlf"
$xlfZhB = "a07768e3-e4d7-49bf-90a4-2429cf1f557e.691f44c004ee09d8f691ab08af1201d70d5f5ea5.zh-cn.xlf"
$xlfDeA = "9e58fb96-cb81-4926-9015-4068e973d72d.81b0d9db6138142e445d0a5a643a729d74d95ba8.de-de.xlf"
$xlfDeB = "a07768e3-e4d7-49bf-90a4-2429cf1f557e.691f44c004ee09d8f691ab08af1201d70d5f5ea5.de-de.xlf"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$include         = "Include"
$zeroDate        = "0001-01-01 00:00:00"
$zhDatetime      = "2016-02-22 17:13:45"
$deDatetime      = "2016-02-22 17:13:56"

# ---------------------------------------------------------------------------
# Sheet "Overview" — simple 3-column file list
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

# Duplicate the last row (".localization-config") twice, pushing it from row 4
# down to row 6 while keeping its original formatting/style untouched.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert(-4121)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert(-4121)

$ws.Range("A4").Value2 = $newMdA
$ws.Range("B4").Value2 = $readyForHandoff
$ws.Range("C4").Value2 = $readyForHandoff

$ws.Range("A5").Value2 = $newMdB
$ws.Range("B5").Value2 = $readyForHandoff
$ws.Range("C5").Value2 = $readyForHandoff

$ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + "f129241a-cb71-4564-963a-2930f3cdfa8b.md"), "", "", "f129241a-cb71-4564-963a-2930f3cdfa8b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + "f13ccd86-f44d-453d-ba24-87b664e3626c.md"), "", "", "f13ccd86-f44d-453d-ba24-87b664e3626c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), ($e2eBase + $newMdA), "", "", $newMdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), ($e2eBase + $newMdB), "", "", $newMdB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), ($commitBase + ".localization-config"), "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" — full handoff/handback tracking table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert(-4121)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert(-4121)

$ws.Range("A4").Value2 = $newMdA
$ws.Range("B4").Value2 = $readyForHandoff
$ws.Range("C4").Value2 = $xlfZhA
$ws.Range("D4").Value2 = $zhDatetime
$ws.Range("G4").Value2 = $zeroDate
$ws.Range("H4").Value2 = $include

$ws.Range("A5").Value2 = $newMdB
$ws.Range("B5").Value2 = $readyForHandoff
$ws.Range("C5").Value2 = $xlfZhB
$ws.Range("D5").Value2 = $zhDatetime
$ws.Range("G5").Value2 = $zeroDate
$ws.Range("H5").Value2 = $include

$ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + "f129241a-cb71-4564-963a-2930f3cdfa8b.md"), "", "", "f129241a-cb71-4564-963a-2930f3cdfa8b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), ($handoffZh + "f129241a-cb71-4564-963a-2930f3cdfa8b.a51f2836575bd58134d2aacbe2bf3f314e0c6e30.zh-cn.xlf"), "", "", "f129241a-cb71-4564-963a-2930f3cdfa8b.a51f2836575bd58134d2aacbe2bf3f314e0c6e30.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + "f13ccd86-f44d-453d-ba24-87b664e3626c.md"), "", "", "f13ccd86-f44d-453d-ba24-87b664e3626c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), ($handoffZh + "f13ccd86-f44d-453d-ba24-87b664e3626c.76ddf37a1611cd8c186558bf124caa745a2a8fa9.zh-cn.xlf"), "", "", "f13ccd86-f44d-453d-ba24-87b664e3626c.76ddf37a1611cd8c186558bf124caa745a2a8fa9.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), ($e2eBase + $newMdA), "", "", $newMdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), ($handoffZh + $xlfZhA), "", "", $xlfZhA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), ($e2eBase + $newMdB), "", "", $newMdB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), ($handoffZh + $xlfZhB), "", "", $xlfZhB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), ($commitBase + ".localization-config"), "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" — full handoff/handback tracking table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert(-4121)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert(-4121)

$ws.Range("A4").Value2 = $newMdA
$ws.Range("B4").Value2 = $readyForHandoff
$ws.Range("C4").Value2 = $xlfDeA
$ws.Range("D4").Value2 = $deDatetime
$ws.Range("G4").Value2 = $zeroDate
$ws.Range("H4").Value2 = $include

$ws.Range("A5").Value2 = $newMdB
$ws.Range("B5").Value2 = $readyForHandoff
$ws.Range("C5").Value2 = $xlfDeB
$ws.Range("D5").Value2 = $deDatetime
$ws.Range("G5").Value2 = $zeroDate
$ws.Range("H5").Value2 = $include

$ws.Hyperlinks.Add($ws.Range("A2"), ($e2eBase + "f129241a-cb71-4564-963a-2930f3cdfa8b.md"), "", "", "f129241a-cb71-4564-963a-2930f3cdfa8b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), ($handoffDe + "f129241a-cb71-4564-963a-2930f3cdfa8b.a51f2836575bd58134d2aacbe2bf3f314e0c6e30.de-de.xlf"), "", "", "f129241a-cb71-4564-963a-2930f3cdfa8b.a51f2836575bd58134d2aacbe2bf3f314e0c6e30.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($e2eBase + "f13ccd86-f44d-453d-ba24-87b664e3626c.md"), "", "", "f13ccd86-f44d-453d-ba24-87b664e3626c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), ($handoffDe + "f13ccd86-f44d-453d-ba24-87b664e3626c.76ddf37a1611cd8c186558bf124caa745a2a8fa9.de-de.xlf"), "", "", "f13ccd86-f44d-453d-ba24-87b664e3626c.76ddf37a1611cd8c186558bf124caa745a2a8fa9.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), ($e2eBase + $newMdA), "", "", $newMdA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), ($handoffDe + $xlfDeA), "", "", $xlfDeA) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), ($e2eBase + $newMdB), "", "", $newMdB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), ($handoffDe + $xlfDeB), "", "", $xlfDeB) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), ($commitBase + ".localization-config"), "", "", ".localization-config") | Out-Null
